$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Non-breaking space used throughout this BOM's order-number strings
$nbsp = [char]0x00A0

# --- 1. Insert 4 blank rows right after "DC Jack" (row 16) for the new
#        component block (LM334 / 3.3 Ohm 0603 / 100 Ohm 0603 / BC807).
#        These land at rows 17-20 (old summary rows 18-19 push to 22-23). ---
$ws.Rows.Item(17).Resize(4).Insert()

# --- 2. Fill "LM334" (row 17). Note: Order-No. was entered before the
#        Manufacturer Part No. for this particular row. ---
$ws.Range("A17").Value = "LM334"
$ws.Range("C17").Value = "Texas Instruments"
$ws.Range("E17").Value = "1014051$($nbsp)- 62"
$ws.Range("D17").Value = "LM334M/NOPB"
$ws.Range("B17").Value = 1
$ws.Range("F17").Value = 0.65
$ws.Range("G17").Value = 0.58
$ws.Range("H17").Formula = "=IF(B`$22*B17>=10,G17,F17)*(B`$22*B17)"

# --- 3. Fill "3.3 Ohm 0603" (row 18) ---
$ws.Range("A18").Value = "3.3 Ohm 0603"
$ws.Range("C18").Value = "Samsung"
$ws.Range("D18").Value = "RC1608F3R3CS"
$ws.Range("E18").Value = "436677$($nbsp)- 62"
$ws.Range("B18").Value = 1
$ws.Range("F18").Value = 0.02
$ws.Range("G18").Value = 0.02
$ws.Range("H18").Formula = "=IF(B`$22*B18>=10,G18,F18)*(B`$22*B18)"

# --- 4. Fill "100 Ohm 0603" (row 19); re-uses the existing "TRU-Components"
#        string, and its part number is numeric. ---
$ws.Range("A19").Value = "100 Ohm 0603"
$ws.Range("C19").Value = "TRU-Components"
$ws.Range("D19").Value = 1585242
$ws.Range("E19").Value = "1585242$($nbsp)- 62"
$ws.Range("B19").Value = 1
$ws.Range("F19").Value = 0.02
$ws.Range("G19").Value = 0.02
$ws.Range("H19").Formula = "=IF(B`$22*B19>=10,G19,F19)*(B`$22*B19)"

# --- 5. Fill "BC807" (row 20) ---
$ws.Range("A20").Value = "BC807"
$ws.Range("C20").Value = "Nexperia"
$ws.Range("D20").Value = "BC807-25,215$($nbsp)"
$ws.Range("E20").Value = "1112884$($nbsp)- 62"
$ws.Range("B20").Value = 1
$ws.Range("F20").Value = 0.2
$ws.Range("G20").Value = 0.19
$ws.Range("H20").Formula = "=IF(B`$22*B20>=10,G20,F20)*(B`$22*B20)"

# --- 6. Insert a new row at 13 for "USB-A" (pushes rows 13-20 -> 14-21,
#        and the summary block 22-23 -> 23-24). Filled in last. ---
$ws.Rows.Item(13).Insert()
$ws.Range("A13").Value = "USB-A"
$ws.Range("C13").Value = "econ-connect"
$ws.Range("D13").Value = "USBBUVA"
$ws.Range("E13").Value = "1311468$($nbsp)- 62"
$ws.Range("B13").Value = 1
$ws.Range("F13").Value = 0.91
$ws.Range("G13").Value = 0.85
$ws.Range("H13").Formula = "=IF(B`$23*B13>=10,G13,F13)*(B`$23*B13)"

# --- 7. Row 22 stays blank except the H-column shared formula still runs
#        through it (Excel drags the formula fill across the blank row). ---
$ws.Range("H22").Formula = "=IF(B`$23*B22>=10,G22,F22)*(B`$23*B22)"

# --- 8. Fix up the summary row (now at 23): the sum range was left at
#        H4:H18 (matching the author's actual edit) rather than growing to
#        cover all the newly added rows. ---
$ws.Range("H23").Formula = "=SUM(H4:H18)"

# --- 9. Restore the selection / active cell as recorded in the saved file ---
$ws.Range("A3:H24").Select()
